$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "mostFrom" text for the bag-for-women row (was "Indonesia")
$ws.Range("D5").Value = "Mainland China"

# Update avgPrice / avgSold values
$ws.Range("B2").Value = 27.7769387755102
$ws.Range("C2").Value = 33773.8163265306

$ws.Range("B3").Value = 20.53625
$ws.Range("C3").Value = 50400.0892857143

$ws.Range("B4").Value = 14.9858823529412
$ws.Range("C4").Value = 38534.7352941177

$ws.Range("B5").Value = 23.2623076923077
$ws.Range("C5").Value = 13229.7692307692

$ws.Range("B6").Value = 17.140625
$ws.Range("C6").Value = 30043.3958333333
